# Updated symbol list on Tue Dec 20 06:51:30 UTC 2022 with GitHub Actions
# Refresh of the cryptocurrency price/volume table on Sheet1.
#
# Column D holds numeric-looking values stored as literal TEXT in the
# source workbook (inlineStr, no leading apostrophe shown to the user).
# Assigning a numeric-looking string straight to .Value would make Excel
# re-interpret it as a real number (and would also force a style change
# the moment a quote-prefix got attached), so each of those cells is
# briefly marked as Text (@ number format), written, then restored to the
# default "Normal" style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Column D price tweaks (rows untouched otherwise) ------------------
$priceUpdates = [ordered]@{
    "D2"  = "248.37"
    "D4"  = "5.343"
    "D5"  = "0.05611"
    "D6"  = "3.407"
    "D7"  = "6.380"
    "D8"  = "0.8158"
    "D9"  = "0.9513"
    "D10" = "0.1419"
    "D11" = "0.07527"
    "D12" = "0.03186"
    "D14" = "0.09310"
    "D15" = "3.571"
    "D16" = "0.001604"
    "D17" = "0.04711"
    "D25" = "0.3298"
    "D26" = "0.1289"
    "D28" = "0.0003000"
    "D40" = "0.03952"
    "D41" = "0.006996"
    "D43" = "0.003030"
    "D44" = "0.008774"
    "D45" = "0.00005727"
    "D47" = "0.0005499"
    "D48" = "0.7799"
    "D49" = "0.1699"
    "D50" = "0.00002100"
}
foreach ($addr in $priceUpdates.Keys) {
    Set-TextValue $addr $priceUpdates[$addr]
}

# --- Rows 18-24: the coin list shifted, "One" now leads off -------------
$rowUpdates = [ordered]@{
    "B18" = "One";        "C18" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one";         "D18" = "0.0005779"; "E18" = "17OneONE"
    "B19" = "TigerCash";  "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";        "D19" = "0.006251";  "E19" = "18TigerCashTCH"
    "B20" = "HotbitToken";"C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";  "D20" = "0.005085";  "E20" = "19HotbitTokenHTB"
    "B21" = "BitKan";     "C21" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";      "D21" = "0.001032";  "E21" = "20BitKanKAN"
    "B22" = "NitroEx";    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";       "D22" = "0.0001500"; "E22" = "21NitroExNTX"
    "B23" = "LEO";        "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";          "D23" = "3.748";     "E23" = "22LEOLEO"
    "B24" = "BTSEToken";  "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";   "D24" = "2.147";     "E24" = "23BTSETokenBTSE"
}
foreach ($addr in $rowUpdates.Keys) {
    $val = $rowUpdates[$addr]
    if ($addr.StartsWith("D")) {
        Set-TextValue $addr $val
    } else {
        $ws.Range($addr).Value = $val
    }
}

# --- Misc text-only tweak -------------------------------------------------
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
